$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/29/2024  Through  2/4/2024"

# --- Pure numeric value updates (crime stat counts / % changes) ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -14.285714285714
$ws.Range("L16").Value = -20
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -80
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = -26.666666666666
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = -50
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 17
$ws.Range("J18").Value = 27
$ws.Range("K18").Value = -37.037037037037
$ws.Range("L18").Value = -39.285714285714
$ws.Range("M18").Value = -5.555555555555
$ws.Range("N18").Value = -87.022900763358
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -11.320754716981
$ws.Range("I19").Value = 58
$ws.Range("J19").Value = 59
$ws.Range("K19").Value = -1.694915254237
$ws.Range("L19").Value = 56.756756756756
$ws.Range("M19").Value = 141.666666666667
$ws.Range("N19").Value = 38.095238095238
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = -44.444444444444
$ws.Range("L20").Value = -37.5
$ws.Range("M20").Value = -37.5
$ws.Range("N20").Value = -93.589743589743
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -18.181818181818
$ws.Range("G21").Value = 98
$ws.Range("H21").Value = -13.265306122449
$ws.Range("I21").Value = 103
$ws.Range("J21").Value = 119
$ws.Range("K21").Value = -13.44537815126
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 41.095890410958
$ws.Range("N21").Value = -69.253731343283
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = -50
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 29.411764705882
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 64
$ws.Range("H24").Value = 14.0625
$ws.Range("I24").Value = 87
$ws.Range("J24").Value = 74
$ws.Range("K24").Value = 17.567567567567
$ws.Range("L24").Value = -1.136363636363
$ws.Range("M24").Value = 102.325581395349
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 100
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 21.052631578947
$ws.Range("I25").Value = 27
$ws.Range("J25").Value = 23
$ws.Range("K25").Value = 17.391304347826
$ws.Range("L25").Value = 22.727272727272
$ws.Range("M25").Value = 68.75
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("K27").Value = -85.714285714285
$ws.Range("L27").Value = -50

# --- Cells changing from a number to the "no data" text markers ("0" / "***.*") ---
# Force text storage via Text number format, then restore the General look by
# copying the format from an existing "no data" text cell (style source: C14).
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cells changing from the "no data" text markers back to real numbers ---
# Restore numeric formatting by copying the format from an existing numeric cell
# (style source: C16 for counts "#,##0", E16 for percentages "#,##0.0").
$ws.Range("E16").Copy()
$ws.Range("L23").PasteSpecial(-4122)
$ws.Range("L23").Value = 0
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("I27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("E16").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("C16").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("E16").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("C16").Copy()
$ws.Range("J30").PasteSpecial(-4122)
$ws.Range("J30").Value = 1
$ws.Range("E16").Copy()
$ws.Range("K30").PasteSpecial(-4122)
$ws.Range("K30").Value = -100
$excel.CutCopyMode = $false
